# edit.ps1 - apply "Sprint 1" deck edits via PowerPoint COM interop
#
# Summary of changes (see task diff):
#   Slide 4 ("Data Dictionary" slide)
#     - Title text "Data Dictionary" -> "Data"
#     - Table row "Rain?" / Data Type cell: "int64" -> "int64 (previously object)"
#     - TextBox 4 ("One concern ... duplicates.") resized + two new paragraphs added
#   Slide 7 ("Next Steps" slide)
#     - Content placeholder gets a new paragraph inserted after the
#       "seasonality of sales." bullet

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 4
# ---------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)

# --- Title: "Data Dictionary" -> "Data" ---
$title4 = $slide4.Shapes.Item(1)
$title4.TextFrame.TextRange.Text = "Data"

# --- Table: row "Rain?", column "Data Type": "int64" -> "int64 (previously object)" ---
$table4 = $slide4.Shapes.Item(2).Table
$table4.Cell(10, 3).Shape.TextFrame.TextRange.Text = "int64 (previously object)"

# --- TextBox 4: update text (drop trailing spaces on line 1, add 2 new lines) and resize ---
$note4 = $slide4.Shapes.Item(3)
$note4.TextFrame.TextRange.Text = "One concern for data quality is that there are 203 duplicates." + [char]13 + "No seasonality feature. " + [char]13 + "  "

# Reposition/resize the box to match the grown text (values picked so the
# saved EMU, which is derived from a single-precision point value, lands
# exactly on the target integers).
$note4.Left = 194.85709461417323
$note4.Top = 474.8984251968504
$note4.Width = 590.5714423228347
$note4.Height = 72.70315260629921

# ---------------------------------------------------------------------
# Slide 7
# ---------------------------------------------------------------------
$slide7 = $p.Slides.Item(7)
$content7 = $slide7.Shapes.Item(3)
$tr7 = $content7.TextFrame.TextRange

$newText7 = "Potentially encoding the different flavors within the data set. " + [char]13 + "Explore feature engineering to create a column associated with the seasonality of sales. " + [char]13 + "Continuing Exploratory Data Analysis within the context of time. " + [char]13 + "Explore different Time Series models to understand which one works best in terms of our proposed solution. "
$tr7.Text = $newText7
